# Kharagauli municipality - "Average monthly remuneration" sheet
# Upgrade the left table with the 2023 column (year header + the three
# figures: total / women / men), matching the layout already used for
# 2014..2022 in columns B..J.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kharagauli")

# Duplicate the formatting of column J (2022) into the new column K (2023)
# so the number format / borders / alignment line up with the rest of the
# table.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the 2023 data for the new column.
$ws.Cells.Item(3, 11).Value = 2023     # header year
$ws.Cells.Item(4, 11).Value = 1609.2   # total average monthly remuneration
$ws.Cells.Item(5, 11).Value = 419.4    # women
$ws.Cells.Item(6, 11).Value = 1824.8   # men
